$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 7) for the "QAP" pipeline, right after the "ANTs" row.
$ws.Range("B7").Value = "QAP"
$ws.Range("C7").Value = 0.055
$ws.Range("D7").Value = 0.000001
$ws.Range("E7").Value = 0.000001
$ws.Range("F7").Value = 20.5
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = "c3.4xlarge"
$ws.Range("I7").Value = 4

# Update the selected/active cell to E12, matching the saved view state.
$ws.Range("E12").Select()
